$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 455.2143
$ws.Range("J92").Value = 527
$ws.Range("L92").Value = 527
$ws.Range("N92").Value = -3023
$ws.Range("H138").Value = 7580666
$ws.Range("I138").Value = 1295.6154
$ws.Range("K138").Value = 3886.8462
$ws.Range("M138").Value = 1253.1538
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 665.7059
$ws.Range("I2").Value = 644.8125
$ws.Range("J2").Value = 1000
$ws.Range("K2").Value = 644.8125
$ws.Range("L2").Value = 1000
$ws.Range("M2").Value = -531.8125
$ws.Range("N2").Value = -1226
$ws.Range("H45").Value = 1956.8096
$ws.Range("I45").Value = 1750.9375
$ws.Range("K45").Value = 1750.9375
$ws.Range("M45").Value = -1373.9375
$ws.Range("H74").Value = 21764102
$ws.Range("I74").Value = 22753306
$ws.Range("K74").Value = 22753306
$ws.Range("M74").Value = -22752432
$ws.Range("H77").Value = 21764102
$ws.Range("I77").Value = 22753306
$ws.Range("K77").Value = 113766530
$ws.Range("M77").Value = -113762162
$ws.Range("H116").Value = 665.7059
$ws.Range("I116").Value = 644.8125
$ws.Range("J116").Value = 1000
$ws.Range("K116").Value = 644.8125
$ws.Range("L116").Value = 1000
$ws.Range("M116").Value = 1649.1875
$ws.Range("N116").Value = -5588
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 665.7059
$ws.Range("I3").Value = 644.8125
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 644.8125
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -530.8125
$ws.Range("N3").Value = -1228
$ws.Range("H20").Value = 2141
$ws.Range("I20").Value = 2150.889
$ws.Range("K20").Value = 2150.889
$ws.Range("M20").Value = -1903.889
$ws.Range("H94").Value = 1592.8125
$ws.Range("I94").Value = 448.2857
$ws.Range("K94").Value = 448.2857
$ws.Range("M94").Value = 2.71429999999998
$ws.Range("H110").Value = 58997.5
$ws.Range("J110").Value = 58997.5
$ws.Range("L110").Value = 58997.5
$ws.Range("N110").Value = -67177.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 10865.538
$ws.Range("I105").Value = 1859.1666
$ws.Range("J105").Value = 18585.285
$ws.Range("K105").Value = 1859.1666
$ws.Range("L105").Value = 18585.285
$ws.Range("M105").Value = -112.1666
$ws.Range("N105").Value = -22079.285
$ws.Range("H110").Value = 37333
$ws.Range("J110").Value = 45999.5
$ws.Range("L110").Value = 45999.5
$ws.Range("N110").Value = -54179.5
$ws.Range("H116").Value = 48664.332
$ws.Range("J116").Value = 48664.332
$ws.Range("L116").Value = 48664.332
$ws.Range("N116").Value = -57842.332
$ws.Range("H122").Value = 1714.6842
$ws.Range("I122").Value = 1463.4375
$ws.Range("J122").Value = 3054.6667
$ws.Range("K122").Value = 4390.3125
$ws.Range("L122").Value = 9164.000100000001
$ws.Range("M122").Value = -1940.3125
$ws.Range("N122").Value = -14064.0001
$ws.Range("H132").Value = 55989.79
$ws.Range("I132").Value = 67478.06
$ws.Range("J132").Value = 5113.143
$ws.Range("K132").Value = 202434.18
$ws.Range("L132").Value = 15339.429
$ws.Range("M132").Value = -199904.18
$ws.Range("N132").Value = -20399.429
$ws.Range("H141").Value = 133075.08
$ws.Range("J141").Value = 133075.08
$ws.Range("L141").Value = 133075.08
$ws.Range("N141").Value = -143435.08
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 4416.6665
$ws.Range("I112").Value = 4416.6665
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 13249.9995
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = -12141.9995
$ws.Range("N112").Value = $null
$ws.Range("H131").Value = 1687.9736
$ws.Range("J131").Value = 1720.3823
$ws.Range("L131").Value = 5161.1469
$ws.Range("N131").Value = -15241.1469
$ws.Range("H133").Value = 4345
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = $null
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = $null
$ws.Range("H80").Value = 3166.8125
$ws.Range("I80").Value = 2953.25
$ws.Range("K80").Value = 2953.25
$ws.Range("M80").Value = -1955.25
$ws.Range("H83").Value = 3166.8125
$ws.Range("I83").Value = 2953.25
$ws.Range("K83").Value = 14766.25
$ws.Range("M83").Value = -9774.25
$ws.Range("H97").Value = 5127.5
$ws.Range("I97").Value = 4855.8
$ws.Range("J97").Value = 5399.2
$ws.Range("K97").Value = 4855.8
$ws.Range("L97").Value = 5399.2
$ws.Range("M97").Value = -4359.8
$ws.Range("N97").Value = -6391.2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3412.0908
$ws.Range("I7").Value = 2646.5
$ws.Range("K7").Value = 2646.5
$ws.Range("M7").Value = -2534.5
$ws.Range("H74").Value = 59996.332
$ws.Range("I74").Value = 61990
$ws.Range("K74").Value = 61990
$ws.Range("M74").Value = -60992
$ws.Range("H75").Value = 46333.332
$ws.Range("J75").Value = 25000
$ws.Range("L75").Value = 25000
$ws.Range("N75").Value = -26872
$ws.Range("H77").Value = 59996.332
$ws.Range("I77").Value = 61990
$ws.Range("K77").Value = 185970
$ws.Range("M77").Value = -180978
$ws.Range("H78").Value = 46333.332
$ws.Range("J78").Value = 25000
$ws.Range("L78").Value = 75000
$ws.Range("N78").Value = -84360
$ws.Range("H122").Value = 3862.9119
$ws.Range("I122").Value = 2359.8235
$ws.Range("K122").Value = 7079.470499999999
$ws.Range("M122").Value = -4629.470499999999
$ws.Range("H126").Value = 3412.0908
$ws.Range("I126").Value = 2646.5
$ws.Range("K126").Value = 7939.5
$ws.Range("M126").Value = -5469.5
$ws.Range("H131").Value = 79677
$ws.Range("J131").Value = 79677
$ws.Range("L131").Value = 79677
$ws.Range("N131").Value = -89757
$ws.Range("H133").Value = 67613
$ws.Range("J133").Value = 67613
$ws.Range("L133").Value = 67613
$ws.Range("N133").Value = -72673
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 80000
$ws.Range("J70").Value = 80000
$ws.Range("L70").Value = 80000
$ws.Range("N70").Value = -80630
$ws.Range("H73").Value = 80000
$ws.Range("J73").Value = 80000
$ws.Range("L73").Value = 80000
$ws.Range("N73").Value = -82184
$ws.Range("H102").Value = 55000
$ws.Range("J102").Value = 50000
$ws.Range("L102").Value = 50000
$ws.Range("N102").Value = -56490
$ws.Range("H122").Value = 49355.76
$ws.Range("I122").Value = 57162.445
$ws.Range("J122").Value = 2515.6667
$ws.Range("K122").Value = 171487.335
$ws.Range("L122").Value = 7547.000100000001
$ws.Range("M122").Value = -169037.335
$ws.Range("N122").Value = -12447.0001
$ws.Range("H132").Value = 2537.818
$ws.Range("I132").Value = 2692.15
$ws.Range("K132").Value = 8076.450000000001
$ws.Range("M132").Value = -5546.450000000001
